# 10.02/2024 - return to host 93
# Update tyre stock sheet: refresh model/size/param rotation and push the
# sales date from 2024-02-29 (45351) to 2024-03-20 (45371); append new
# models/rows so the sheet grows from A1:J30 to A1:J34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'BEL-262'
$ws.Range('B2').Value = '205/55R16'
$ws.Range('C2').Value = 'легк сер б/к'
$ws.Range('E2').Value = '205/55R16'
$ws.Range('F2').Value = 'BEL-262'
$ws.Range('G2').Value = 'легк, сер, б/к'
$ws.Range("I2").Value = 45371
$ws.Range('A3').Value = 'BEL-317'
$ws.Range('C3').Value = 'легк сер б/к'
$ws.Range('E3').Value = '205/55R16'
$ws.Range('F3').Value = 'BEL-317'
$ws.Range('G3').Value = 'легк, сер, б/к'
$ws.Range("I3").Value = 45371
$ws.Range('A4').Value = 'BEL-317S'
$ws.Range('C4').Value = 'сер ошип'
$ws.Range('F4').Value = 'BEL-317S'
$ws.Range('G4').Value = 'сер, ошип'
$ws.Range("I4").Value = 45371
$ws.Range('A5').Value = 'BEL-1001'
$ws.Range('B5').Value = '235/75R15'
$ws.Range('C5').Value = 'легк сер'
$ws.Range('E5').Value = '235/75R15'
$ws.Range('F5').Value = 'BEL-1001'
$ws.Range('G5').Value = 'легк, сер'
$ws.Range("I5").Value = 45371
$ws.Range('A6').Value = 'BEL-1002'
$ws.Range('B6').Value = '155/65R13'
$ws.Range('E6').Value = '155/65R13'
$ws.Range('F6').Value = 'BEL-1002'
$ws.Range('G6').Value = 'легк, сер'
$ws.Range("I6").Value = 45371
$ws.Range('A7').Value = 'BEL-1004'
$ws.Range('B7').Value = '205/55R16'
$ws.Range('E7').Value = '205/55R16'
$ws.Range('F7').Value = 'BEL-1004'
$ws.Range("I7").Value = 45371
$ws.Range('A8').Value = 'BEL-1005'
$ws.Range('B8').Value = '225/50R17'
$ws.Range('E8').Value = '205/55R16'
$ws.Range('F8').Value = 'BEL-1004'
$ws.Range("I8").Value = 45371
$ws.Range('A9').Value = 'Бел-202'
$ws.Range('B9').Value = '24.00R35'
$ws.Range('C9').Value = '210B Type H сер'
$ws.Range('E9').Value = '225/50R17'
$ws.Range('F9').Value = 'BEL-1005'
$ws.Range("I9").Value = 45371
$ws.Range('A10').Value = 'BEL-248'
$ws.Range('B10').Value = '14.00R20'
$ws.Range('C10').Value = 'груз сер б/к'
$ws.Range('E10').Value = '24.00R35'
$ws.Range('F10').Value = 'Бел-202'
$ws.Range('G10').Value = '210B, Type, H, сер'
$ws.Range("I10").Value = 45371
$ws.Range('A11').Value = 'Бел-103'
$ws.Range('B11').Value = '175/70R13'
$ws.Range('C11').Value = 'легк сер б/к'
$ws.Range('E11').Value = '14.00R20'
$ws.Range('F11').Value = 'BEL-248'
$ws.Range('G11').Value = 'груз, сер, б/к'
$ws.Range("I11").Value = 45371
$ws.Range('A12').Value = 'Бел-100'
$ws.Range('B12').Value = '175/70R13'
$ws.Range('C12').Value = 'легк сер б/к'
$ws.Range('E12').Value = '14.00R20'
$ws.Range('F12').Value = 'BEL-248'
$ws.Range('G12').Value = 'груз, сер, б/к'
$ws.Range("I12").Value = 45371
$ws.Range('A13').Value = 'Ф-35-1'
$ws.Range('B13').Value = '11.2-20'
$ws.Range('C13').Value = '8 сх сер'
$ws.Range('E13').Value = '14.00R20'
$ws.Range('F13').Value = 'BEL-248'
$ws.Range('G13').Value = 'груз, сер, б/к'
$ws.Range("I13").Value = 45371
$ws.Range('A14').Value = 'Бел-119'
$ws.Range('B14').Value = '195/65R15'
$ws.Range('C14').Value = 'легк сер'
$ws.Range('E14').Value = '14.00R20'
$ws.Range('F14').Value = 'BEL-248'
$ws.Range('G14').Value = 'груз, сер, б/к'
$ws.Range("I14").Value = 45371
$ws.Range('A15').Value = 'Бел-1149'
$ws.Range('B15').Value = '195/65R15'
$ws.Range('C15').Value = 'легк сер'
$ws.Range('E15').Value = '175/70R13'
$ws.Range('F15').Value = 'Бел-103'
$ws.Range('G15').Value = 'легк, сер, б/к'
$ws.Range("I15").Value = 45371
$ws.Range('A16').Value = 'Бел-777'
$ws.Range('B16').Value = '210/80R16'
$ws.Range('C16').Value = 'легк сер'
$ws.Range('E16').Value = '175/70R13'
$ws.Range('F16').Value = 'Бел-100'
$ws.Range('G16').Value = 'легк, сер, б/к'
$ws.Range("I16").Value = 45371
$ws.Range('A17').Value = 'Бел-1000'
$ws.Range('B17').Value = '215/65R16C'
$ws.Range('C17').Value = 'легк сер'
$ws.Range('E17').Value = '195/65R15'
$ws.Range('F17').Value = 'Бел-119'
$ws.Range('G17').Value = 'легк, сер'
$ws.Range("I17").Value = 45371
$ws.Range('A18').Value = 'Бел-1001'
$ws.Range('B18').Value = '205/55R16'
$ws.Range('E18').Value = '195/65R15'
$ws.Range('F18').Value = 'Бел-119'
$ws.Range('G18').Value = 'легк, сер'
$ws.Range("I18").Value = 45371
$ws.Range('A19').Value = 'Бел-1005'
$ws.Range('B19').Value = '225/50R17'
$ws.Range('E19').Value = '195/65R15'
$ws.Range('F19').Value = 'Бел-119'
$ws.Range('G19').Value = 'легк, сер'
$ws.Range("I19").Value = 45371
$ws.Range('A20').Value = 'BEL-734'
$ws.Range('B20').Value = '205/70R14'
$ws.Range('C20').Value = 'легк сер б/к'
$ws.Range('E20').Value = '195/65R15'
$ws.Range('F20').Value = 'Бел-1149'
$ws.Range('G20').Value = 'легк, сер'
$ws.Range("I20").Value = 45371
$ws.Range('A21').Value = 'BEL-261'
$ws.Range('B21').Value = '195/65R15'
$ws.Range('C21').Value = 'легк сер б/к'
$ws.Range('E21').Value = '210/80R16'
$ws.Range('F21').Value = 'Бел-777'
$ws.Range('G21').Value = 'легк, сер'
$ws.Range("I21").Value = 45371
$ws.Range('A22').Value = 'BEL-337'
$ws.Range('B22').Value = '195/65R15'
$ws.Range('C22').Value = 'легк сер б/к'
$ws.Range('E22').Value = '215/65R16C'
$ws.Range('F22').Value = 'Бел-1000'
$ws.Range('G22').Value = 'легк, сер'
$ws.Range("I22").Value = 45371
$ws.Range('A23').Value = 'BEL-337S'
$ws.Range('B23').Value = '195/65R15'
$ws.Range('C23').Value = 'сер ошип'
$ws.Range('E23').Value = '205/55R16'
$ws.Range('F23').Value = 'Бел-1001'
$ws.Range('G23').Value = 'легк, сер'
$ws.Range("I23").Value = 45371
$ws.Range('A24').Value = 'BEL-705'
$ws.Range('B24').Value = '195/65R15'
$ws.Range('C24').Value = 'легк сер б/к'
$ws.Range('E24').Value = '225/50R17'
$ws.Range('F24').Value = 'Бел-1005'
$ws.Range("I24").Value = 45371
$ws.Range('A25').Value = 'Бел-188'
$ws.Range('B25').Value = '175/70R13'
$ws.Range('C25').Value = 'легк сер'
$ws.Range('E25').Value = '205/70R14'
$ws.Range('F25').Value = 'BEL-734'
$ws.Range('G25').Value = 'легк, сер, б/к'
$ws.Range("I25").Value = 45371
$ws.Range('F26').Value = 'BEL-261'
$ws.Range('G26').Value = 'легк, сер, б/к'
$ws.Range("I26").Value = 45371
$ws.Range('E27').Value = '195/65R15'
$ws.Range('F27').Value = 'BEL-337'
$ws.Range('G27').Value = 'легк, сер, б/к'
$ws.Range("I27").Value = 45371
$ws.Range('E28').Value = '195/65R15'
$ws.Range('F28').Value = 'BEL-337S'
$ws.Range('G28').Value = 'сер, ошип'
$ws.Range("I28").Value = 45371
$ws.Range('E29').Value = '195/65R15'
$ws.Range('F29').Value = 'BEL-705'
$ws.Range('G29').Value = 'легк, сер, б/к'
$ws.Range("I29").Value = 45371
$ws.Range('E30').Value = '205/55R16'
$ws.Range('F30').Value = 'BEL-262'
$ws.Range('G30').Value = 'легк, сер, б/к'
$ws.Range("I30").Value = 45371
$ws.Range('E31').Value = '205/55R16'
$ws.Range('F31').Value = 'BEL-317'
$ws.Range('G31').Value = 'легк, сер, б/к'
$ws.Range("H31").Value = 2
$ws.Range("I31").Value = 45371
$ws.Range('J31').Value = 'нет данных'
$ws.Range('E32').Value = '175/70R13'
$ws.Range('F32').Value = 'Бел-103'
$ws.Range('G32').Value = 'легк, сер, б/к'
$ws.Range("H32").Value = 2
$ws.Range("I32").Value = 45371
$ws.Range('J32').Value = 'нет данных'
$ws.Range('E33').Value = '175/70R13'
$ws.Range('F33').Value = 'Бел-100'
$ws.Range('G33').Value = 'легк, сер, б/к'
$ws.Range("H33").Value = 2
$ws.Range("I33").Value = 45371
$ws.Range('J33').Value = 'нет данных'
$ws.Range('E34').Value = '175/70R13'
$ws.Range('F34').Value = 'Бел-188'
$ws.Range('G34').Value = 'легк, сер'
$ws.Range("H34").Value = 2
$ws.Range("I34").Value = 45371
$ws.Range('J34').Value = 'нет данных'

$ws.Range("I23").NumberFormat = "yyyy-mm-dd"
$ws.Range("I24").NumberFormat = "yyyy-mm-dd"
$ws.Range("I25").NumberFormat = "yyyy-mm-dd"
$ws.Range("I31").NumberFormat = "yyyy-mm-dd"
$ws.Range("I32").NumberFormat = "yyyy-mm-dd"
$ws.Range("I33").NumberFormat = "yyyy-mm-dd"
$ws.Range("I34").NumberFormat = "yyyy-mm-dd"
